# Auto-generated edit script applying the Aegis_Profits.xlsx diff
# Updates per-row market board / profit figures across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 355193.3
$ws.Range("I33").Value = 467.3125
$ws.Range("K33").Value = 467.3125
$ws.Range("M33").Value = -238.3125
$ws.Range("H34").Value = 2934.5557
$ws.Range("I34").Value = 2934.5557
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2934.5557
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2731.5557
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 2934.5557
$ws.Range("I36").Value = 2934.5557
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2934.5557
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2219.5557
$ws.Range("N36").ClearContents()
$ws.Range("H61").Value = 3722750
$ws.Range("I61").Value = 4166841.5
$ws.Range("K61").Value = 12500524.5
$ws.Range("M61").Value = -12500352.5
$ws.Range("H107").Value = 649.62964
$ws.Range("I107").Value = 668.26086
$ws.Range("K107").Value = 668.26086
$ws.Range("M107").Value = 1251.73914
$ws.Range("H137").Value = 1504.4348
$ws.Range("I137").Value = 1362
$ws.Range("J137").Value = 1908
$ws.Range("K137").Value = 4086
$ws.Range("L137").Value = 5724
$ws.Range("M137").Value = -1536
$ws.Range("N137").Value = -10824

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 63886.188
$ws.Range("I2").Value = 1512.7858
$ws.Range("J2").Value = 500500
$ws.Range("K2").Value = 1512.7858
$ws.Range("L2").Value = 500500
$ws.Range("M2").Value = -1399.7858
$ws.Range("N2").Value = -500726
$ws.Range("H74").Value = 1130.0555
$ws.Range("I74").Value = 1090.9231
$ws.Range("J74").Value = 1231.8
$ws.Range("K74").Value = 1090.9231
$ws.Range("L74").Value = 1231.8
$ws.Range("M74").Value = -216.9231
$ws.Range("N74").Value = -2979.8
$ws.Range("H77").Value = 1130.0555
$ws.Range("I77").Value = 1090.9231
$ws.Range("J77").Value = 1231.8
$ws.Range("K77").Value = 5454.6155
$ws.Range("L77").Value = 6159
$ws.Range("M77").Value = -1086.6155
$ws.Range("N77").Value = -14895
$ws.Range("H95").Value = 30736
$ws.Range("J95").Value = 30736
$ws.Range("L95").Value = 30736
$ws.Range("N95").Value = -36228
$ws.Range("H116").Value = 63886.188
$ws.Range("I116").Value = 1512.7858
$ws.Range("J116").Value = 500500
$ws.Range("K116").Value = 1512.7858
$ws.Range("L116").Value = 500500
$ws.Range("M116").Value = 781.2141999999999
$ws.Range("N116").Value = -505088
$ws.Range("H132").Value = 13207.8125
$ws.Range("I132").Value = 15381.718
$ws.Range("J132").Value = 3787.5557
$ws.Range("K132").Value = 46145.154
$ws.Range("L132").Value = 11362.6671
$ws.Range("M132").Value = -43615.154
$ws.Range("N132").Value = -16422.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 63886.188
$ws.Range("I3").Value = 1512.7858
$ws.Range("J3").Value = 500500
$ws.Range("K3").Value = 1512.7858
$ws.Range("L3").Value = 500500
$ws.Range("M3").Value = -1398.7858
$ws.Range("N3").Value = -500728
$ws.Range("H26").Value = 27257
$ws.Range("I26").Value = 12885.5
$ws.Range("J26").Value = 56000
$ws.Range("K26").Value = 12885.5
$ws.Range("L26").Value = 56000
$ws.Range("M26").Value = -12593.5
$ws.Range("N26").Value = -56584
$ws.Range("H134").Value = 18187.605
$ws.Range("I134").Value = 21268.223
$ws.Range("J134").Value = 4324.8335
$ws.Range("K134").Value = 63804.66900000001
$ws.Range("L134").Value = 12974.5005
$ws.Range("M134").Value = -61269.66900000001
$ws.Range("N134").Value = -18044.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 29000
$ws.Range("J29").Value = 29000
$ws.Range("L29").Value = 29000
$ws.Range("N29").Value = -29586
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H132").Value = 83337064
$ws.Range("I132").Value = 83337520
$ws.Range("J132").Value = 83336150
$ws.Range("K132").Value = 250012560
$ws.Range("L132").Value = 250008450
$ws.Range("M132").Value = -250010030
$ws.Range("N132").Value = -250013510

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1203.0161
$ws.Range("J5").Value = 1286.9778
$ws.Range("L5").Value = 3860.9334
$ws.Range("N5").Value = -4084.9334
$ws.Range("H113").Value = 856.9666999999999
$ws.Range("I113").Value = 1179.3334
$ws.Range("J113").Value = 534.6
$ws.Range("K113").Value = 3538.0002
$ws.Range("L113").Value = 1603.8
$ws.Range("M113").Value = -1368.0002
$ws.Range("N113").Value = -5943.8
$ws.Range("H131").Value = 828.1
$ws.Range("J131").Value = 861.07526
$ws.Range("L131").Value = 2583.22578
$ws.Range("N131").Value = -12663.22578
$ws.Range("H135").Value = 1203.0161
$ws.Range("J135").Value = 1286.9778
$ws.Range("L135").Value = 11582.8002
$ws.Range("N135").Value = -16652.8002
$ws.Range("H141").Value = 11491
$ws.Range("I141").Value = 13115
$ws.Range("K141").Value = 39345
$ws.Range("M141").Value = -34165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3204.6
$ws.Range("I122").Value = 2503.75
$ws.Range("J122").Value = 6008
$ws.Range("K122").Value = 7511.25
$ws.Range("L122").Value = 18024
$ws.Range("M122").Value = -5061.25
$ws.Range("N122").Value = -22924
$ws.Range("H132").Value = 3851.95
$ws.Range("I132").Value = 3292.2
$ws.Range("J132").Value = 4411.7
$ws.Range("K132").Value = 9876.599999999999
$ws.Range("L132").Value = 13235.1
$ws.Range("M132").Value = -7346.599999999999
$ws.Range("N132").Value = -18295.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3160
$ws.Range("I7").Value = 1688
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 1688
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -1576
$ws.Range("N7").Value = -5224
$ws.Range("H93").Value = 3547.4443
$ws.Range("I93").Value = 3657.1667
$ws.Range("J93").Value = 3328
$ws.Range("K93").Value = 3657.1667
$ws.Range("L93").Value = 3328
$ws.Range("M93").Value = -2409.1667
$ws.Range("N93").Value = -5824
$ws.Range("H126").Value = 3160
$ws.Range("I126").Value = 1688
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5064
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -2594
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 5738.706
$ws.Range("I132").Value = 4826.357
$ws.Range("J132").Value = 9996.333000000001
$ws.Range("K132").Value = 14479.071
$ws.Range("L132").Value = 29988.999
$ws.Range("M132").Value = -11949.071
$ws.Range("N132").Value = -35048.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 250012300
$ws.Range("H42").Value = 29008.166
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 29008.166
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 29008.166
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -29764.166
$ws.Range("H43").Value = 16400
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 16400
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 16400
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -16698
$ws.Range("H132").Value = 2961.0557
$ws.Range("J132").Value = 2816
$ws.Range("L132").Value = 8448
$ws.Range("N132").Value = -13508
